$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-21 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-22 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("490÷2=245, 0", $true, $false, $false, $false, $false, $true, 1, $false, "854÷6=142, 2", 2) | Out-Null
$d.Content.Find.Execute("941÷4=235, 1", $true, $false, $false, $false, $false, $true, 1, $false, "544÷6=90, 4", 2) | Out-Null
$d.Content.Find.Execute("434÷8=54, 2", $true, $false, $false, $false, $false, $true, 1, $false, "213÷8=26, 5", 2) | Out-Null
$d.Content.Find.Execute("250÷2=125, 0", $true, $false, $false, $false, $false, $true, 1, $false, "349÷4=87, 1", 2) | Out-Null
$d.Content.Find.Execute("297÷8=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "162÷7=23, 1", 2) | Out-Null
$d.Content.Find.Execute("760÷8=95, 0", $true, $false, $false, $false, $false, $true, 1, $false, "482÷9=53, 5", 2) | Out-Null
$d.Content.Find.Execute("785÷7=112, 1", $true, $false, $false, $false, $false, $true, 1, $false, "360÷7=51, 3", 2) | Out-Null
$d.Content.Find.Execute("918÷8=114, 6", $true, $false, $false, $false, $false, $true, 1, $false, "273÷6=45, 3", 2) | Out-Null
$d.Content.Find.Execute("542÷8=67, 6", $true, $false, $false, $false, $false, $true, 1, $false, "624÷7=89, 1", 2) | Out-Null
$d.Content.Find.Execute("347÷9=38, 5", $true, $false, $false, $false, $false, $true, 1, $false, "308÷3=102, 2", 2) | Out-Null
$d.Content.Find.Execute("755÷5=151, 0", $true, $false, $false, $false, $false, $true, 1, $false, "763÷8=95, 3", 2) | Out-Null
$d.Content.Find.Execute("163÷6=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "898÷9=99, 7", 2) | Out-Null
$d.Content.Find.Execute("732÷4=183, 0", $true, $false, $false, $false, $false, $true, 1, $false, "157÷6=26, 1", 2) | Out-Null
$d.Content.Find.Execute("632÷9=70, 2", $true, $false, $false, $false, $false, $true, 1, $false, "163÷4=40, 3", 2) | Out-Null
$d.Content.Find.Execute("275÷7=39, 2", $true, $false, $false, $false, $false, $true, 1, $false, "417÷9=46, 3", 2) | Out-Null
$d.Content.Find.Execute("819÷2=409, 1", $true, $false, $false, $false, $false, $true, 1, $false, "502÷3=167, 1", 2) | Out-Null
$d.Content.Find.Execute("755÷3=251, 2", $true, $false, $false, $false, $false, $true, 1, $false, "956÷8=119, 4", 2) | Out-Null
$d.Content.Find.Execute("688÷6=114, 4", $true, $false, $false, $false, $false, $true, 1, $false, "186÷6=31, 0", 2) | Out-Null
$d.Content.Find.Execute("103÷8=12, 7", $true, $false, $false, $false, $false, $true, 1, $false, "996÷2=498, 0", 2) | Out-Null
$d.Content.Find.Execute("475÷2=237, 1", $true, $false, $false, $false, $false, $true, 1, $false, "916÷8=114, 4", 2) | Out-Null
$d.Content.Find.Execute("668÷4=167, 0", $true, $false, $false, $false, $false, $true, 1, $false, "447÷6=74, 3", 2) | Out-Null
$d.Content.Find.Execute("371÷5=74, 1", $true, $false, $false, $false, $false, $true, 1, $false, "604÷2=302, 0", 2) | Out-Null
$d.Content.Find.Execute("893÷2=446, 1", $true, $false, $false, $false, $false, $true, 1, $false, "756÷8=94, 4", 2) | Out-Null
$d.Content.Find.Execute("707÷4=176, 3", $true, $false, $false, $false, $false, $true, 1, $false, "925÷7=132, 1", 2) | Out-Null
$d.Content.Find.Execute("270÷5=54, 0", $true, $false, $false, $false, $false, $true, 1, $false, "467÷5=93, 2", 2) | Out-Null

Write-Output "Replacements applied."
